# "Generate Report for Handoff"
# Refreshes the localization-status report: the old run's generated-file
# UUID (020d2a01-66cb-4750-9580-cb37a40022e6) is replaced everywhere by the
# new run's UUID (9a9450c5-c191-4572-8715-063a4d1820e2), and the handoff /
# handback timestamps are bumped to the new run's times. The hyperlink
# *targets* (they still point at the same already-published commit URL) are
# left untouched - only the visible display text of each hyperlink changes.

$wb = $excel.ActiveWorkbook

$oldId = "020d2a01-66cb-4750-9580-cb37a40022e6"
$newId = "9a9450c5-c191-4572-8715-063a4d1820e2"

$hyperlinkTarget = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/49ff447e6f96eb4a3c79baa484fc6debc2bbd073/e2e/$oldId.md"

function Update-HyperlinkDisplay {
    param($range, $address, $display)

    # Re-point the cell's hyperlink without touching its Address (target
    # URL / relationship id): writing straight to .TextToDisplay or
    # .Address on an existing Hyperlinks.Item leaves the stale link behind
    # and appends a duplicate, so instead the link on this range is removed
    # and a fresh one is added back with the same address and new display
    # text.
    $range.Hyperlinks.Delete()
    $range.Worksheet.Hyperlinks.Add($range, $address, "", "", $display) | Out-Null
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A2").Value = "$newId.md"
$wsOverview.Range("B2").Value = "e2e\$newId.md"
Update-HyperlinkDisplay $wsOverview.Range("B2") $hyperlinkTarget "e2e\$newId.md"
$wsOverview.Range("G2").Value = "2016-09-02 05:02:23"

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A2").Value = "$newId.md"
Update-HyperlinkDisplay $wsZhCn.Range("A2") $hyperlinkTarget "$newId.md"
$wsZhCn.Range("G2").Value = "$newId.8cf3e283a5c6ba46ce8bfeea422d9646d4919405.zh-cn.xlf"
$wsZhCn.Range("H2").Value = "2016-09-02 05:02:19"

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A2").Value = "$newId.md"
Update-HyperlinkDisplay $wsDeDe.Range("A2") $hyperlinkTarget "$newId.md"
$wsDeDe.Range("G2").Value = "$newId.8cf3e283a5c6ba46ce8bfeea422d9646d4919405.de-de.xlf"
$wsDeDe.Range("H2").Value = "2016-09-02 05:02:23"
